# Slide 5 ("There's much more to Apache Isis"), Content Placeholder 2,
# bullet list under "Other features": rename the 2nd sub-bullet from
# "Home page" to "View models" (keep its existing bold/lang formatting).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$para = $shape.TextFrame.TextRange.Paragraphs(2, 1)
$para.Text = "View models"
